$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write AT1 header first (claims the "RC" shared string)
$ws.Range("AT1").Value = "RC"

# Step 2: rename AU1 header
$ws.Range("AU1").Value = "RC values"

# Step 3: delete redundant defined names
$toDelete = @(14,15,16,17,18,19,20,21)
foreach ($idx in $toDelete) {
    $name = "_xlchart.v1.$idx"
    $wb.Names.Item($name).Delete()
}
